$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.293.22"
$ws.Range("E2").Value = "  -1.53%  "

$ws.Range("D3").Value = "2.985.56"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.25%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("D9").Value = "2.984.68"
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D16").Value = "3.486.25"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").Value = "61.339.82"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").Value = "2.987.25"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.683"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.87%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.87%  "

$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.85%  "

$ws.Range("D35").Value = "0.0₃0821"
$ws.Range("E35").Value = "  +3.85%  "

$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.40%  "

$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "396.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0352"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.269"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("D47").Value = "2.683.89"
$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("E48").Value = "  +2.36%  "

$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "
